$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new "ServicesCounterParty" sheet right after
#    "RelatedCounterParty" (and before "Settlement").
# ------------------------------------------------------------------
$related = $wb.Worksheets.Item("RelatedCounterParty")
$ws = $wb.Worksheets.Add($null, $related)
$ws.Name = "ServicesCounterParty"

# ------------------------------------------------------------------
# 2. Header row (row 1)
# ------------------------------------------------------------------
$ws.Range("A1").Value = "TestCaseID"
$ws.Range("B1").Value = "Action"
$ws.Range("C1").Value = "ServiceRef"
$ws.Range("D1").Value = "Counterparty"
$ws.Range("E1").Value = "BusinessUnit"
$ws.Range("F1").Value = "ServiceName"
$ws.Range("G1").Value = "ServiceType"
$ws.Range("H1").Value = "ServiceTier"
$ws.Range("I1").Value = "FeeRate"
$ws.Range("J1").Value = "FeeType"
$ws.Range("K1").Value = "ChargeRate"
$ws.Range("L1").Value = "ChargeType"
$ws.Range("M1").Value = "ServiceEntity"
$ws.Range("N1").Value = "ServiceStatus"
$ws.Range("O1").Value = "ServiceSignerName"
$ws.Range("P1").Value = "ServiceStartDate"
$ws.Range("Q1").Value = "ServiceEndDate"
$ws.Range("R1").Value = "ServiceContactPerson"
$ws.Range("S1").Value = "ServiceContactPersonEmail"
$ws.Range("T1").Value = "ServiceContactPersonMobile"
$ws.Range("U1").Value = "ExpectedTradingVolume"
$ws.Range("V1").Value = "MiningServiceTransferFee"
$ws.Range("W1").Value = "MiningServiceHashrate"
$ws.Range("X1").Value = "MiningServicePower"
$ws.Range("Y1").Value = "MiningServicePowerCharge"
$ws.Range("Z1").Value = "ReferralPerson"
$ws.Range("AA1").Value = "ServiceRemarks"
$ws.Range("AB1").Value = "Messages"
$ws.Range("AC1").Value = "SkipAtStepNum"
$ws.Range("AD1").Value = "SIT"

# ------------------------------------------------------------------
# 3. Data rows (rows 2-7) - column A first, then column B, which is
#    the order the original authoring tool entered the data in (and
#    which drives shared-string allocation order on save).
# ------------------------------------------------------------------
$ws.Range("A2").Value = "QA"
$ws.Range("A3").Value = "QA_TestCase_Auto_Optimus_2_3_1"
$ws.Range("A4").Value = "QA_TestCase_Auto_Optimus_2_3_2"
$ws.Range("A5").Value = "QA_TestCase_Auto_Optimus_2_3_3"
$ws.Range("A6").Value = "QA_TestCase_Auto_Optimus_2_3_4"
$ws.Range("A7").Value = "QA_TestCase_Auto_Optimus_2_3_5"

$ws.Range("B3").Value = "Create"
$ws.Range("B4").Value = "Edit"
$ws.Range("B5").Value = "Delete"
$ws.Range("B6").Value = "Search"
$ws.Range("B7").Value = "Download Btn"

# ------------------------------------------------------------------
# 4. Header styling - reuse the look & feel of the header cells on
#    "RelatedCounterParty" (same workbook theme / border / bold style)
# ------------------------------------------------------------------
$related.Range("A1").Copy()
$ws.Range("A1").PasteSpecial(-4122) | Out-Null

$related.Range("B1").Copy()
$ws.Range("B1").PasteSpecial(-4122) | Out-Null

$related.Range("AL1").Copy()
$ws.Range("AB1").PasteSpecial(-4122) | Out-Null

$related.Range("AM1").Copy()
$ws.Range("AC1").PasteSpecial(-4122) | Out-Null

$related.Range("AN1").Copy()
$ws.Range("AD1").PasteSpecial(-4122) | Out-Null

$ws.Application.CutCopyMode = $false

# ------------------------------------------------------------------
# 5. Column widths - best-fit the used columns (A:AD)
# ------------------------------------------------------------------
for ($i = 1; $i -le 30; $i++) {
    $ws.Columns.Item($i).EntireColumn.AutoFit() | Out-Null
}

# ------------------------------------------------------------------
# 6. Sheet view / selection bookkeeping
# ------------------------------------------------------------------
# RelatedCounterParty is no longer the active tab; its old cell
# selection becomes a full first-row selection instead.
$related.Rows.Item(1).Select()

# The new sheet becomes the active tab, with B8 selected (matching
# the author's last editing position).
$ws.Select()
$ws.Range("B8").Select()
